# A new price record is inserted as row 114 (pushing the existing rows
# 114-214 down to 115-215). All of the "key" columns (Mercado ID,
# Mercado, Region, Codreg, Tipo, Producto ID, Producto, Categoria ID,
# Categoria, Variedad) are identical for every row in this sheet, so we
# simply insert a full row and fill in the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 114; everything currently at/under
# row 114 shifts down by one (row 114 -> 115, ..., row 214 -> 215).
$ws.Rows.Item(114).Insert()

# Populate the new row 114 with the new price record.
$ws.Cells.Item(114, 1).Value  = 3
$ws.Cells.Item(114, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(114, 3).Value  = "Coquimbo"
$ws.Cells.Item(114, 4).Value  = 44651
$ws.Cells.Item(114, 5).Value  = 5
$ws.Cells.Item(114, 6).Value  = "Fruta"
$ws.Cells.Item(114, 7).Value  = 100101
$ws.Cells.Item(114, 8).Value  = "Berries"
$ws.Cells.Item(114, 9).Value  = 100101001
$ws.Cells.Item(114, 10).Value = "Arándano (blue)"
$ws.Cells.Item(114, 11).Value = "Sin especificar"
$ws.Cells.Item(114, 12).Value = "Primera"
$ws.Cells.Item(114, 13).Value = 30
$ws.Cells.Item(114, 14).Value = 4500
$ws.Cells.Item(114, 15).Value = 4500
$ws.Cells.Item(114, 16).Value = 4500
$ws.Cells.Item(114, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(114, 18).Value = "Provincia de Linares"
$ws.Cells.Item(114, 19).Value = 2250
$ws.Cells.Item(114, 20).Value = 2
